# Auto-generated Excel COM-interop script
# Applies numeric cell updates across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 555.5294
$ws.Range("I9").Value = 585.25
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 585.25
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = -416.25
$ws.Range("N9").Value = -418

$ws.Range("H33").Value = 51504.25
$ws.Range("I33").Value = 71793.07000000001
$ws.Range("J33").Value = 4163.6665
$ws.Range("K33").Value = 71793.07000000001
$ws.Range("L33").Value = 4163.6665
$ws.Range("M33").Value = -71564.07000000001
$ws.Range("N33").Value = -4621.6665

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").Value = $null

$ws.Range("H129").Value = 2254.52
$ws.Range("J129").Value = 1116.7819
$ws.Range("L129").Value = 3350.3457
$ws.Range("N129").Value = -13350.3457

$ws.Range("H132").Value = 7582326
$ws.Range("I132").Value = 9623475
$ws.Range("J132").Value = 914.8570999999999
$ws.Range("K132").Value = 28870425
$ws.Range("L132").Value = 2744.5713
$ws.Range("M132").Value = -28867895
$ws.Range("N132").Value = -7804.5713

$ws.Range("H135").Value = 987.5741
$ws.Range("I135").Value = 679.449
$ws.Range("K135").Value = 6115.040999999999
$ws.Range("M135").Value = -3580.040999999999

$ws.Range("H137").Value = 973.2166999999999
$ws.Range("I137").Value = 953.875
$ws.Range("K137").Value = 2861.625
$ws.Range("M137").Value = -311.625

$ws.Range("H138").Value = 1719.2051
$ws.Range("I138").Value = 1276.0322
$ws.Range("J138").Value = 3436.5
$ws.Range("K138").Value = 3828.0966
$ws.Range("L138").Value = 10309.5
$ws.Range("M138").Value = 1311.9034
$ws.Range("N138").Value = -20589.5

$ws.Range("H141").Value = 1333.6094
$ws.Range("I141").Value = 1252.4354
$ws.Range("J141").Value = 3850
$ws.Range("K141").Value = 3757.3062
$ws.Range("L141").Value = 11550
$ws.Range("M141").Value = 1422.6938
$ws.Range("N141").Value = -21910

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 99.5
$ws.Range("J4").Value = 99
$ws.Range("L4").Value = 99
$ws.Range("N4").Value = -331

$ws.Range("H32").Value = 23477.646
$ws.Range("I32").Value = 5183.375
$ws.Range("J32").Value = 155196.4
$ws.Range("K32").Value = 5183.375
$ws.Range("L32").Value = 155196.4
$ws.Range("M32").Value = -4896.375
$ws.Range("N32").Value = -155770.4

$ws.Range("H37").Value = 11633.333
$ws.Range("I37").Value = 4900
$ws.Range("K37").Value = 4900
$ws.Range("M37").Value = -4627

$ws.Range("H44").Value = 14985
$ws.Range("J44").Value = 14985
$ws.Range("L44").Value = 14985
$ws.Range("N44").Value = -15961

$ws.Range("H55").Value = 13433.333
$ws.Range("I55").Value = 8000
$ws.Range("J55").Value = 14985.714
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 14985.714
$ws.Range("M55").Value = -7685
$ws.Range("N55").Value = -15615.714

$ws.Range("H74").Value = 545.5
$ws.Range("I74").Value = 510.55
$ws.Range("K74").Value = 510.55
$ws.Range("M74").Value = 363.45

$ws.Range("H77").Value = 545.5
$ws.Range("I77").Value = 510.55
$ws.Range("K77").Value = 2552.75
$ws.Range("M77").Value = 1815.25

$ws.Range("H80").Value = 27431.5
$ws.Range("J80").Value = 27431.5
$ws.Range("L80").Value = 27431.5
$ws.Range("N80").Value = -29427.5

$ws.Range("H83").Value = 27431.5
$ws.Range("J83").Value = 27431.5
$ws.Range("L83").Value = 82294.5
$ws.Range("N83").Value = -92278.5

$ws.Range("H139").Value = 67857.5
$ws.Range("J139").Value = 67857.5
$ws.Range("L139").Value = 67857.5
$ws.Range("N139").Value = -78137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1237
$ws.Range("J16").Value = 2747.5
$ws.Range("L16").Value = 2747.5
$ws.Range("N16").Value = -3321.5

$ws.Range("H20").Value = 42936.715
$ws.Range("J20").Value = 42936.715
$ws.Range("L20").Value = 42936.715
$ws.Range("N20").Value = -43408.715

$ws.Range("H30").Value = 42936.715
$ws.Range("J30").Value = 42936.715
$ws.Range("L30").Value = 42936.715
$ws.Range("N30").Value = -43118.715

$ws.Range("H31").Value = 23844.385
$ws.Range("I31").Value = 1327.6316
$ws.Range("K31").Value = 1327.6316
$ws.Range("M31").Value = -1032.6316

$ws.Range("H34").Value = 23844.385
$ws.Range("I34").Value = 1327.6316
$ws.Range("K34").Value = 1327.6316
$ws.Range("M34").Value = -1125.6316

$ws.Range("H113").Value = 1237
$ws.Range("J113").Value = 2747.5
$ws.Range("L113").Value = 2747.5
$ws.Range("N113").Value = -7087.5

$ws.Range("H128").Value = 42936.715
$ws.Range("J128").Value = 42936.715
$ws.Range("L128").Value = 42936.715
$ws.Range("N128").Value = -52896.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 720
$ws.Range("I34").Value = 553
$ws.Range("J34").Value = 794.2222
$ws.Range("K34").Value = 1659
$ws.Range("L34").Value = 2382.6666
$ws.Range("M34").Value = -1575
$ws.Range("N34").Value = -2550.6666

$ws.Range("H80").Value = 12221.444
$ws.Range("J80").Value = 12221.444
$ws.Range("L80").Value = 36664.33199999999
$ws.Range("N80").Value = -38536.33199999999

$ws.Range("H83").Value = 12221.444
$ws.Range("J83").Value = 12221.444
$ws.Range("L83").Value = 109992.996
$ws.Range("N83").Value = -119352.996

$ws.Range("H131").Value = 8643.026
$ws.Range("J131").Value = 8691.6
$ws.Range("L131").Value = 26074.8
$ws.Range("N131").Value = -36154.8

$ws.Range("H134").Value = 3137.4285
$ws.Range("I134").Value = 2802
$ws.Range("J134").Value = 3271.6
$ws.Range("K134").Value = 8406
$ws.Range("L134").Value = 9814.799999999999
$ws.Range("M134").Value = -3336
$ws.Range("N134").Value = -19954.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5850143
$ws.Range("I11").Value = 5850143
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5850143
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -5850004
$ws.Range("N11").Value = $null

$ws.Range("H113").Value = 1633.8096
$ws.Range("J113").Value = 1490.5385
$ws.Range("L113").Value = 1490.5385
$ws.Range("N113").Value = -5830.538500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3150867
$ws.Range("I16").Value = 3818848
$ws.Range("K16").Value = 3818848
$ws.Range("M16").Value = -3818678

$ws.Range("H20").Value = 48337.332
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 70006
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 70006
$ws.Range("M20").Value = -4774
$ws.Range("N20").Value = -70458

$ws.Range("H24").Value = 24900
$ws.Range("J24").Value = 24900
$ws.Range("L24").Value = 24900
$ws.Range("N24").Value = -25586

$ws.Range("H25").Value = 200037310
$ws.Range("I25").Value = 500008500
$ws.Range("J25").Value = 56533.332
$ws.Range("K25").Value = 500008500
$ws.Range("L25").Value = 56533.332
$ws.Range("M25").Value = -500008270
$ws.Range("N25").Value = -56993.332

$ws.Range("H100").Value = 2186.4443
$ws.Range("I100").Value = 1860
$ws.Range("J100").Value = 2594.5
$ws.Range("K100").Value = 1860
$ws.Range("L100").Value = 2594.5
$ws.Range("M100").Value = -1319
$ws.Range("N100").Value = -3676.5

$ws.Range("H132").Value = 2384.077
$ws.Range("I132").Value = 2308.2444
$ws.Range("J132").Value = 2871.5715
$ws.Range("K132").Value = 6924.733200000001
$ws.Range("L132").Value = 8614.7145
$ws.Range("M132").Value = -4394.733200000001
$ws.Range("N132").Value = -13674.7145

$ws.Range("H136").Value = 1146.0392
$ws.Range("I136").Value = 921.36365
$ws.Range("J136").Value = 2558.2856
$ws.Range("K136").Value = 2764.09095
$ws.Range("L136").Value = 7674.8568
$ws.Range("M136").Value = -214.0909499999998
$ws.Range("N136").Value = -12774.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 10751
$ws.Range("I20").Value = 7710
$ws.Range("J20").Value = 11088.889
$ws.Range("K20").Value = 7710
$ws.Range("L20").Value = 11088.889
$ws.Range("M20").Value = -7470
$ws.Range("N20").Value = -11568.889

$ws.Range("H132").Value = 1765.6666
$ws.Range("I132").Value = 1671
$ws.Range("J132").Value = 2295.8
$ws.Range("K132").Value = 5013
$ws.Range("L132").Value = 6887.400000000001
$ws.Range("M132").Value = -2483
$ws.Range("N132").Value = -11947.4

$ws.Range("H136").Value = 511.93442
$ws.Range("I136").Value = 348.7143
$ws.Range("J136").Value = 1178.4166
$ws.Range("K136").Value = 1046.1429
$ws.Range("L136").Value = 3535.2498
$ws.Range("M136").Value = 1503.8571
$ws.Range("N136").Value = -8635.2498
